# echoMRI editions in timedura column and CD1 exploration updates
#
# 1) Reformat the TimeDateDura strings in column G (rows 2-31) from
#      "HH:MM:SS; 1 Aug 2025; N; ems"
#    to
#      "HH:MM:SS Aug 1, 2025; N; ems"
#    The G23 cell is intentionally updated last so that the shared-string
#    table ends up ordered the same way the source workbook was produced
#    (all the other rows first, top to bottom, then row 23 last).
#
# 2) Cosmetic workbook/worksheet view updates: zoom to 140%, select G24,
#    and widen column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the TimeDateDura (column G) text values -------------------

$rowOrder = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,24,25,26,27,28,29,30,31,23)

foreach ($r in $rowOrder) {
    $cell = $ws.Cells.Item($r, 7)
    $old = $cell.Value2
    if ($old -match '^(\d{1,2}:\d{2}:\d{2}); 1 Aug 2025(.*)$') {
        $cell.Value = "$($matches[1]) Aug 1, 2025$($matches[2])"
    }
}

# --- 2) View / layout tweaks ----------------------------------------------

$win = $excel.ActiveWindow
$win.Zoom = 140

$ws.Range("G24").Select()

$ws.Columns.Item(7).ColumnWidth = 24.71
